$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.036344
$ws.Range("I2").Value = 0.00971307879529387
$ws.Range("J2").Value = 0.009713078795293874
$ws.Range("M2").Value = 9.358310000000001
$ws.Range("N2").Value = 28.07493
$ws.Range("O2").Value = 0.1719944618809179
$ws.Range("P2").Value = 0.1719944618809179
$ws.Range("Q2").Value = 0.34011841864
$ws.Range("R2").Value = 3.06106576776
$ws.Range("S2").Value = 0.001670595760603524
$ws.Range("T2").Value = 0.001670595760603524
$ws.Range("G3").Value = 0.036344
$ws.Range("I3").Value = 0.00971307879529387
$ws.Range("J3").Value = 0.009713078795293874
$ws.Range("O3").Value = 0.4866586274141183
$ws.Range("P3").Value = 0.4866586274141184
$ws.Range("Q3").Value = 0.9623656538906668
$ws.Range("R3").Value = 8.661290885016001
$ws.Range("S3").Value = 0.004726953594482892
$ws.Range("T3").Value = 0.004726953594482895
$ws.Range("G4").Value = 0.036344
$ws.Range("I4").Value = 0.00971307879529387
$ws.Range("J4").Value = 0.009713078795293874
$ws.Range("M4").Value = 5.152806666666667
$ws.Range("N4").Value = 15.45842
$ws.Range("O4").Value = 0.09470237786627494
$ws.Range("P4").Value = 0.09470237786627496
$ws.Range("Q4").Value = 0.1872736054933334
$ws.Range("R4").Value = 1.68546244944
$ws.Range("S4").Value = 0.0009198516583168228
$ws.Range("T4").Value = 0.0009198516583168232
$ws.Range("G5").Value = 0.036344
$ws.Range("I5").Value = 0.00971307879529387
$ws.Range("J5").Value = 0.009713078795293874
$ws.Range("M5").Value = 7.974813333333334
$ws.Range("N5").Value = 23.92444
$ws.Range("O5").Value = 0.146567460136225
$ws.Range("P5").Value = 0.146567460136225
$ws.Range("Q5").Value = 0.2898366157866667
$ws.Range("R5").Value = 2.60852954208
$ws.Range("S5").Value = 0.001423621289129247
$ws.Range("T5").Value = 0.001423621289129247
$ws.Range("G6").Value = 0.036344
$ws.Range("I6").Value = 0.00971307879529387
$ws.Range("J6").Value = 0.009713078795293874
$ws.Range("M6").Value = 5.445246666666667
$ws.Range("N6").Value = 16.33574
$ws.Range("O6").Value = 0.1000770727024639
$ws.Range("P6").Value = 0.1000770727024639
$ws.Range("Q6").Value = 0.1979020448533333
$ws.Range("R6").Value = 1.78111840368
$ws.Range("S6").Value = 0.0009720564927613854
$ws.Range("T6").Value = 0.0009720564927613858
$ws.Range("I7").Value = 0.8635746806069907
$ws.Range("J7").Value = 0.8635746806069909
$ws.Range("M7").Value = 9.358310000000001
$ws.Range("N7").Value = 28.07493
$ws.Range("O7").Value = 0.1719944618809179
$ws.Range("P7").Value = 0.1719944618809179
$ws.Range("Q7").Value = 30.23939792271667
$ws.Range("R7").Value = 272.15458130445
$ws.Range("S7").Value = 0.1485300624849849
$ws.Range("T7").Value = 0.148530062484985
$ws.Range("I8").Value = 0.8635746806069907
$ws.Range("J8").Value = 0.8635746806069909
$ws.Range("O8").Value = 0.4866586274141183
$ws.Range("P8").Value = 0.4866586274141184
$ws.Range("S8").Value = 0.4202660687337837
$ws.Range("T8").Value = 0.4202660687337839
$ws.Range("I9").Value = 0.8635746806069907
$ws.Range("J9").Value = 0.8635746806069909
$ws.Range("M9").Value = 5.152806666666667
$ws.Range("N9").Value = 15.45842
$ws.Range("O9").Value = 0.09470237786627494
$ws.Range("P9").Value = 0.09470237786627496
$ws.Range("Q9").Value = 16.65020406592222
$ws.Range("R9").Value = 149.8518365933
$ws.Range("S9").Value = 0.08178257571859093
$ws.Range("T9").Value = 0.08178257571859096
$ws.Range("I10").Value = 0.8635746806069907
$ws.Range("J10").Value = 0.8635746806069909
$ws.Range("M10").Value = 7.974813333333334
$ws.Range("N10").Value = 23.92444
$ws.Range("O10").Value = 0.146567460136225
$ws.Range("P10").Value = 0.146567460136225
$ws.Range("Q10").Value = 25.76892128451111
$ws.Range("R10").Value = 231.9202915606
$ws.Range("S10").Value = 0.1265719475745183
$ws.Range("T10").Value = 0.1265719475745184
$ws.Range("I11").Value = 0.8635746806069907
$ws.Range("J11").Value = 0.8635746806069909
$ws.Range("M11").Value = 5.445246666666667
$ws.Range("N11").Value = 16.33574
$ws.Range("O11").Value = 0.1000770727024639
$ws.Range("P11").Value = 0.1000770727024639
$ws.Range("Q11").Value = 17.59516202612222
$ws.Range("R11").Value = 158.3564582351
$ws.Range("S11").Value = 0.08642402609511288
$ws.Range("T11").Value = 0.0864240260951129
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.2040783333333333
$ws.Range("H12").Value = 0.612235
$ws.Range("I12").Value = 0.05454074763589353
$ws.Range("J12").Value = 0.05454074763589354
$ws.Range("M12").Value = 9.358310000000001
$ws.Range("N12").Value = 28.07493
$ws.Range("O12").Value = 0.1719944618809179
$ws.Range("P12").Value = 0.1719944618809179
$ws.Range("Q12").Value = 1.909828307616667
$ws.Range("R12").Value = 17.18845476855
$ws.Range("S12").Value = 0.00938070654021845
$ws.Range("T12").Value = 0.009380706540218456
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.2040783333333333
$ws.Range("H13").Value = 0.612235
$ws.Range("I13").Value = 0.05454074763589353
$ws.Range("J13").Value = 0.05454074763589354
$ws.Range("O13").Value = 0.4866586274141183
$ws.Range("P13").Value = 0.4866586274141184
$ws.Range("Q13").Value = 5.403862500089445
$ws.Range("R13").Value = 48.634762500805
$ws.Range("S13").Value = 0.02654272538262376
$ws.Range("T13").Value = 0.02654272538262377
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.2040783333333333
$ws.Range("H14").Value = 0.612235
$ws.Range("I14").Value = 0.05454074763589353
$ws.Range("J14").Value = 0.05454074763589354
$ws.Range("M14").Value = 5.152806666666667
$ws.Range("N14").Value = 15.45842
$ws.Range("O14").Value = 0.09470237786627494
$ws.Range("P14").Value = 0.09470237786627496
$ws.Range("Q14").Value = 1.051576196522222
$ws.Range("R14").Value = 9.4641857687
$ws.Range("S14").Value = 0.00516513849172353
$ws.Range("T14").Value = 0.005165138491723532
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.2040783333333333
$ws.Range("H15").Value = 0.612235
$ws.Range("I15").Value = 0.05454074763589353
$ws.Range("J15").Value = 0.05454074763589354
$ws.Range("M15").Value = 7.974813333333334
$ws.Range("N15").Value = 23.92444
$ws.Range("O15").Value = 0.146567460136225
$ws.Range("P15").Value = 0.146567460136225
$ws.Range("Q15").Value = 1.627486613711111
$ws.Range("R15").Value = 14.6473795234
$ws.Range("S15").Value = 0.007993898854923732
$ws.Range("T15").Value = 0.007993898854923735
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.2040783333333333
$ws.Range("H16").Value = 0.612235
$ws.Range("I16").Value = 0.05454074763589353
$ws.Range("J16").Value = 0.05454074763589354
$ws.Range("M16").Value = 5.445246666666667
$ws.Range("N16").Value = 16.33574
$ws.Range("O16").Value = 0.1000770727024639
$ws.Range("P16").Value = 0.1000770727024639
$ws.Range("Q16").Value = 1.111256864322222
$ws.Range("R16").Value = 10.0013117789
$ws.Range("S16").Value = 0.005458278366404053
$ws.Range("T16").Value = 0.005458278366404056
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1463416666666667
$ws.Range("H17").Value = 0.439025
$ws.Range("I17").Value = 0.03911039344507936
$ws.Range("J17").Value = 0.03911039344507936
$ws.Range("M17").Value = 9.358310000000001
$ws.Range("N17").Value = 28.07493
$ws.Range("O17").Value = 0.1719944618809179
$ws.Range("P17").Value = 0.1719944618809179
$ws.Range("Q17").Value = 1.369510682583334
$ws.Range("R17").Value = 12.32559614325
$ws.Range("S17").Value = 0.006726771074537402
$ws.Range("T17").Value = 0.006726771074537404
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.1463416666666667
$ws.Range("H18").Value = 0.439025
$ws.Range("I18").Value = 0.03911039344507936
$ws.Range("J18").Value = 0.03911039344507936
$ws.Range("O18").Value = 0.4866586274141183
$ws.Range("P18").Value = 0.4866586274141184
$ws.Range("Q18").Value = 3.875032845397222
$ws.Range("R18").Value = 34.875295608575
$ws.Range("S18").Value = 0.01903341039160845
$ws.Range("T18").Value = 0.01903341039160845
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.1463416666666667
$ws.Range("H19").Value = 0.439025
$ws.Range("I19").Value = 0.03911039344507936
$ws.Range("J19").Value = 0.03911039344507936
$ws.Range("M19").Value = 5.152806666666667
$ws.Range("N19").Value = 15.45842
$ws.Range("O19").Value = 0.09470237786627494
$ws.Range("P19").Value = 0.09470237786627496
$ws.Range("Q19").Value = 0.7540703156111112
$ws.Range("R19").Value = 6.7866328405
$ws.Range("S19").Value = 0.003703847258534588
$ws.Range("T19").Value = 0.003703847258534589
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.1463416666666667
$ws.Range("H20").Value = 0.439025
$ws.Range("I20").Value = 0.03911039344507936
$ws.Range("J20").Value = 0.03911039344507936
$ws.Range("M20").Value = 7.974813333333334
$ws.Range("N20").Value = 23.92444
$ws.Range("O20").Value = 0.146567460136225
$ws.Range("P20").Value = 0.146567460136225
$ws.Range("Q20").Value = 1.167047474555556
$ws.Range("R20").Value = 10.503427271
$ws.Range("S20").Value = 0.005732311032173743
$ws.Range("T20").Value = 0.005732311032173746
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.1463416666666667
$ws.Range("H21").Value = 0.439025
$ws.Range("I21").Value = 0.03911039344507936
$ws.Range("J21").Value = 0.03911039344507936
$ws.Range("M21").Value = 5.445246666666667
$ws.Range("N21").Value = 16.33574
$ws.Range("O21").Value = 0.1000770727024639
$ws.Range("P21").Value = 0.1000770727024639
$ws.Range("Q21").Value = 0.7968664726111112
$ws.Range("R21").Value = 7.1717982535
$ws.Range("S21").Value = 0.003914053688225175
$ws.Range("T21").Value = 0.003914053688225176
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0.6666666666666666
$ws.Range("G22").Value = 0.1237066666666667
$ws.Range("H22").Value = 0.37112
$ws.Range("I22").Value = 0.03306109951674244
$ws.Range("J22").Value = 0.03306109951674245
$ws.Range("M22").Value = 9.358310000000001
$ws.Range("N22").Value = 28.07493
$ws.Range("O22").Value = 0.1719944618809179
$ws.Range("P22").Value = 0.1719944618809179
$ws.Range("Q22").Value = 1.157685335733334
$ws.Range("R22").Value = 10.4191680216
$ws.Range("S22").Value = 0.00568632602057359
$ws.Range("T22").Value = 0.005686326020573593
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0.6666666666666666
$ws.Range("G23").Value = 0.1237066666666667
$ws.Range("H23").Value = 0.37112
$ws.Range("I23").Value = 0.03306109951674244
$ws.Range("J23").Value = 0.03306109951674245
$ws.Range("O23").Value = 0.4866586274141183
$ws.Range("P23").Value = 0.4866586274141184
$ws.Range("Q23").Value = 3.275672660062222
$ws.Range("R23").Value = 29.48105394056
$ws.Range("S23").Value = 0.01608946931161944
$ws.Range("T23").Value = 0.01608946931161945
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0.6666666666666666
$ws.Range("G24").Value = 0.1237066666666667
$ws.Range("H24").Value = 0.37112
$ws.Range("I24").Value = 0.03306109951674244
$ws.Range("J24").Value = 0.03306109951674245
$ws.Range("M24").Value = 5.152806666666667
$ws.Range("N24").Value = 15.45842
$ws.Range("O24").Value = 0.09470237786627494
$ws.Range("P24").Value = 0.09470237786627496
$ws.Range("Q24").Value = 0.6374365367111111
$ws.Range("R24").Value = 5.7369288304
$ws.Range("S24").Value = 0.003130964739109062
$ws.Range("T24").Value = 0.003130964739109064
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 0.6666666666666666
$ws.Range("G25").Value = 0.1237066666666667
$ws.Range("H25").Value = 0.37112
$ws.Range("I25").Value = 0.03306109951674244
$ws.Range("J25").Value = 0.03306109951674245
$ws.Range("M25").Value = 7.974813333333334
$ws.Range("N25").Value = 23.92444
$ws.Range("O25").Value = 0.146567460136225
$ws.Range("P25").Value = 0.146567460136225
$ws.Range("Q25").Value = 0.9865375747555556
$ws.Range("R25").Value = 8.8788381728
$ws.Range("S25").Value = 0.004845681385479915
$ws.Range("T25").Value = 0.004845681385479917
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 0.6666666666666666
$ws.Range("G26").Value = 0.1237066666666667
$ws.Range("H26").Value = 0.37112
$ws.Range("I26").Value = 0.03306109951674244
$ws.Range("J26").Value = 0.03306109951674245
$ws.Range("M26").Value = 5.445246666666667
$ws.Range("N26").Value = 16.33574
$ws.Range("O26").Value = 0.1000770727024639
$ws.Range("P26").Value = 0.1000770727024639
$ws.Range("Q26").Value = 0.6736133143111112
$ws.Range("R26").Value = 6.0625198288
$ws.Range("S26").Value = 0.003308658059960428
$ws.Range("T26").Value = 0.00330865805996043
